# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to match the latest scrape (gh-pages output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet ---
$wsExhibition.Range("F2").Value = 1884
$wsExhibition.Range("F3").Value = 497
$wsExhibition.Range("F4").Value = 169
$wsExhibition.Range("F5").Value = 176
$wsExhibition.Range("F6").Value = 2617
$wsExhibition.Range("F7").Value = 173
$wsExhibition.Range("F8").Value = 93
$wsExhibition.Range("F10").Value = 1543
$wsExhibition.Range("F11").Value = 537
$wsExhibition.Range("F13").Value = 336
$wsExhibition.Range("F14").Value = 232
$wsExhibition.Range("F15").Value = 23
$wsExhibition.Range("F18").Value = 220
$wsExhibition.Range("F21").Value = 184
$wsExhibition.Range("F22").Value = 61
$wsExhibition.Range("F23").Value = 1677
$wsExhibition.Range("F24").Value = 35
$wsExhibition.Range("F25").Value = 411
$wsExhibition.Range("F26").Value = 15
$wsExhibition.Range("F28").Value = 210
$wsExhibition.Range("F29").Value = 303
$wsExhibition.Range("F30").Value = 424

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F2").Value = 1884
$wsAll.Range("F4").Value = 497
$wsAll.Range("F5").Value = 169
$wsAll.Range("F6").Value = 176
$wsAll.Range("F7").Value = 2617
$wsAll.Range("F8").Value = 173
$wsAll.Range("F9").Value = 93
$wsAll.Range("F11").Value = 1543
$wsAll.Range("F12").Value = 537
$wsAll.Range("F14").Value = 336
$wsAll.Range("F15").Value = 232
$wsAll.Range("F16").Value = 23
$wsAll.Range("F19").Value = 220
$wsAll.Range("F22").Value = 184
$wsAll.Range("F23").Value = 61
$wsAll.Range("F24").Value = 1677
$wsAll.Range("F25").Value = 35
$wsAll.Range("F26").Value = 411
$wsAll.Range("F27").Value = 15
$wsAll.Range("F29").Value = 210
$wsAll.Range("F30").Value = 303
$wsAll.Range("F31").Value = 424

